# pedalboard-display-bom.xlsx update
# - Split the "BoM" sheet's J1/J2 connector row into a fitted J1 (LED-IN) row
#   that stays on BoM, and a new "DNF" (Did Not Fit) sheet holding J2 (LED-OUT).
# - Update the summary counts on BoM (Component Groups/Fitted/Total).
# - New "DNF" sheet is inserted between "BoM" and "Colors", carrying the same
#   header/branding block, picture and frozen header row as BoM.

$wb = $excel.ActiveWorkbook
$bom = $wb.Worksheets.Item("BoM")

# ---------------------------------------------------------------------------
# 1. Update BoM summary numbers
# ---------------------------------------------------------------------------
$bom.Range("F2").Value2 = 7                          # Component Groups: 6 -> 7
$bom.Range("F4").Value2 = "201 (199 SMD/ 2 THT)"      # Fitted Components:
$bom.Range("F6").Value2 = 201                         # Total Components: 202 -> 201

# ---------------------------------------------------------------------------
# 2. Split BoM row 11 (J1 J2 / Conn_01x03_Pin) into just J1 / LED-IN, qty 1
# ---------------------------------------------------------------------------
$bom.Range("D11").Value2 = "J1"
$bom.Range("E11").Value2 = "LED-IN"
$bom.Range("G11").Value2 = 1

# ---------------------------------------------------------------------------
# 3. Create the new "DNF" sheet, positioned right after "BoM"
#    (copying BoM gives us identical header block, column widths, frozen
#    pane and picture "for free")
# ---------------------------------------------------------------------------
$bom.Copy([System.Reflection.Missing]::Value, $bom)
$dnf = $wb.Worksheets.Item(2)
$dnf.Name = "DNF"

# Row 9 on DNF should look like the old "J1 J2" row (same styling as BoM's
# row 11), but describing J2 / LED-OUT as "Did Not Fit".
$dnf.Range("A11:J11").Copy()
$dnf.Range("A9:J9").PasteSpecial(-4122)   # xlPasteFormats

$dnf.Range("A9").Value2 = 1
$dnf.Range("B9").Value2 = ""
$dnf.Range("C9").Value2 = "Conn_01x03_Pin"
$dnf.Range("D9").Value2 = "J2"
$dnf.Range("E9").Value2 = "LED-OUT"
$dnf.Range("F9").Value2 = "JST_PH_B3B-PH-SM4-TB_1x03-1MP_P2.00mm_Vertical"
$dnf.Range("G9").Value2 = 1
$dnf.Range("H9").Value2 = " (DNF)"
$dnf.Range("I9").Value2 = "~"
$dnf.Range("J9").Value2 = ""

# Remove the remaining BoM rows (10-14) that were copied along with the sheet
$dnf.Rows("10:14").Delete()

# Row 9 shouldn't keep BoM's tall "105pt" row height
$dnf.Rows("9:9").AutoFit()

# Print titles for the new sheet (row 9 is the header row, like on BoM)
$dnf.PageSetup.PrintTitleRows = "`$9:`$9"

# Re-activate BoM as the selected tab (sheet copy leaves the new sheet active)
$bom.Activate()
